# Apply the "Add tests from another source and make change one case" edit:
#  1. Update the date-range header text (07/06/2019 -> 28/06/2019).
#  2. Append 15 new exchange-rate rows (111-125, 10/06/2019 .. 28/06/2019)
#     in the same style as the existing data rows.
#  3. Move the active selection to F16 (matches the saved sheetView state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header text update (shared string used by cell A2).
$ws.Range("A2").Value = "טווח תאריכים : 02/01/2019 - 28/06/2019"

# 2. Append the new rows, carrying over the date/rate number formatting
#    from the last existing row (110) before writing the new values.
$lastRow = 110
$firstNew = $lastRow + 1

$dates = @(43626,43627,43628,43629,43630,43633,43634,43635,43636,43637,43640,43641,43642,43643,43644)
$rates = @(3.585,3.581,3.5819999999999999,3.5920000000000001,3.6,3.61,3.6120000000000001,3.609,3.5790000000000002,3.5939999999999999,3.6040000000000001,3.6019999999999999,3.5910000000000002,3.5819999999999999,3.5659999999999998)

$lastNew = $firstNew + $dates.Length - 1

$ws.Range("A$lastRow`:B$lastRow").Copy()
$ws.Range("A$firstNew`:B$lastNew").PasteSpecial(-4122)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $firstNew + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $rates[$i]
}

# 3. Update the selected / active cell.
$ws.Range("F16").Select()
